# Update the log write mode: refresh run_time (C) and the max_er /
# per-iteration error columns (E:Y) with the latest simulated run data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.5515866279602051
$rowVals = @(763.422004028851, 0.0291592111127273, 0.02431361535913722, 0.02177951082962003, 0.02003573029606406, 0.01946759317727915, 0.01842830622997898, 0.01811816810957158, 0.01746033365562736, 0.01675998033754384, 0.01614792266184621, 0.01583686515905356, 0.01562031831073966, 0.01547484634008248, 0.01523928980627402, 0.01523928980627402, 0.01521880520459229, 0.01507589137739381, 0.0149531288808537, 0.01494433044889471, 0.01488152054637136)
$arr = New-Object 'object[,]' 1,$rowVals.Length
for ($i = 0; $i -lt $rowVals.Length; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("E2:Y2").Value = $arr

$ws.Range("C3").Value = 0.5312278270721436
$rowVals = @(758.9649216285943, 0.02919310067863152, 0.02395173217772089, 0.02240587448097198, 0.02065471345684584, 0.01921349161781889, 0.01831714892507718, 0.01788906728454151, 0.01755175954823061, 0.01707398797071064, 0.01650204385099904, 0.01588895695031006, 0.01558452566071163, 0.01558452566071163, 0.01535079459020347, 0.01527589677106527, 0.01508091276864604, 0.01505110543095962, 0.01496826062795641, 0.01488548934235925, 0.01479463784851061)
$arr = New-Object 'object[,]' 1,$rowVals.Length
for ($i = 0; $i -lt $rowVals.Length; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("E3:Y3").Value = $arr

$ws.Range("C4").Value = 0.6627647876739502
$rowVals = @(772.3454484991344, 0.02850826916610752, 0.02437627267256313, 0.02124452028173654, 0.01970331012014087, 0.01917412434197139, 0.01827847498123716, 0.01732504200754465, 0.01715313501785395, 0.01661122532776104, 0.01638774192204171, 0.01638774192204171, 0.0159877809852551, 0.01585234186679805, 0.01554308533675904, 0.01548961947156114, 0.01525764549510535, 0.01525764549510535, 0.01519141226536524, 0.0150554668323418, 0.0150554668323418)
$arr = New-Object 'object[,]' 1,$rowVals.Length
for ($i = 0; $i -lt $rowVals.Length; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("E4:Y4").Value = $arr

$ws.Range("C5").Value = 0.5480189323425293
$rowVals = @(768.7557478054841, 0.02858914194274109, 0.02289459797593941, 0.02059383023386789, 0.01995797896127936, 0.01864693947562354, 0.01832024079235138, 0.01727626067605043, 0.01718903956280533, 0.01637833514881277, 0.01637833514881277, 0.01568435669629203, 0.0156835389704169, 0.01546074179356897, 0.01516924774570739, 0.0150472627187991, 0.0150472627187991, 0.0150472627187991, 0.01503054058263039, 0.01501289857020675, 0.01498549215995095)
$arr = New-Object 'object[,]' 1,$rowVals.Length
for ($i = 0; $i -lt $rowVals.Length; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("E5:Y5").Value = $arr

$ws.Range("C6").Value = 0.6875057220458984
$rowVals = @(773.2677809646284, 0.02873372457987006, 0.02325534950990853, 0.02193742082708855, 0.01931502422015857, 0.01847940424367292, 0.01786723515483411, 0.01721644798277264, 0.01721644798277264, 0.01679395622489825, 0.0164316081489612, 0.01641766800801105, 0.01589992298999089, 0.0156416689924988, 0.0156416689924988, 0.01555317079354341, 0.01529150713697932, 0.01529150713697932, 0.01516563825096894, 0.01510733573563288, 0.0150734460227023)
$arr = New-Object 'object[,]' 1,$rowVals.Length
for ($i = 0; $i -lt $rowVals.Length; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("E6:Y6").Value = $arr

$ws.Range("C7").Value = 0.7499716281890869
$rowVals = @(757.4912447191273, 0.02880877708122558, 0.02498516550167713, 0.02136163391252212, 0.02049073270265218, 0.01966586239899683, 0.01840371723420983, 0.01781602001346617, 0.01741029191311592, 0.01631339834849248, 0.01601670924476266, 0.01593298684517554, 0.01566631544272797, 0.01544279152291834, 0.01527313443532203, 0.01489171819480735, 0.01489171819480735, 0.01481385779960585, 0.01477480423638571, 0.01477480423638571, 0.01476591120310189)
$arr = New-Object 'object[,]' 1,$rowVals.Length
for ($i = 0; $i -lt $rowVals.Length; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("E7:Y7").Value = $arr

$ws.Range("C8").Value = 0.6718742847442627
$rowVals = @(784.5242615559964, 0.02890530376814575, 0.02428823862717855, 0.02206179127280525, 0.02124969692788641, 0.01944666747266235, 0.0184352012782944, 0.01792144057672153, 0.01778415966491387, 0.01725290766047362, 0.01713085394779034, 0.01654021255728863, 0.01634700363040347, 0.01573363286279296, 0.01573363286279296, 0.01564359667377769, 0.01555207004715836, 0.01546397119322473, 0.01538007651708735, 0.01536376620065324, 0.01529287059563346)
$arr = New-Object 'object[,]' 1,$rowVals.Length
for ($i = 0; $i -lt $rowVals.Length; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("E8:Y8").Value = $arr

$ws.Range("C9").Value = 0.6406235694885254
$rowVals = @(769.194913445408, 0.02880802070235775, 0.02386165180329171, 0.02142225736210604, 0.02019862242082252, 0.01910216296161166, 0.01842003499199273, 0.01797263552390944, 0.01666937600362633, 0.01666937600362633, 0.01660970949427341, 0.01617918952150247, 0.01591713569783143, 0.01567296497253808, 0.01537975044052354, 0.01525442955620374, 0.01518995618218199, 0.01511507542041535, 0.01505081112943678, 0.01499494661371312, 0.01499405289367267)
$arr = New-Object 'object[,]' 1,$rowVals.Length
for ($i = 0; $i -lt $rowVals.Length; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("E9:Y9").Value = $arr

$ws.Range("C10").Value = 0.5937929153442383
$rowVals = @(763.2117246748076, 0.02852373839707052, 0.02304180885986327, 0.02107036949474379, 0.01970042405779887, 0.0187826303641657, 0.01804015584283191, 0.01677156001290929, 0.01616285869850623, 0.01598107378101337, 0.01554133689408702, 0.01554133689408702, 0.01535960166036978, 0.01535960166036978, 0.01531667352002395, 0.01522062358735035, 0.01514386462187435, 0.0150706206537953, 0.01503679995319366, 0.01494994099992449, 0.01487742153362198)
$arr = New-Object 'object[,]' 1,$rowVals.Length
for ($i = 0; $i -lt $rowVals.Length; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("E10:Y10").Value = $arr

$ws.Range("C11").Value = 0.6093349456787109
$rowVals = @(751.6617062712121, 0.02911580777740332, 0.02482051390245419, 0.02169419381170808, 0.01993946896506335, 0.01899268016579646, 0.01852992231920022, 0.01666341676627872, 0.01608812272556948, 0.01606382401946655, 0.01585480815690737, 0.01559528081540926, 0.01537341703699672, 0.01531918964780433, 0.01520815314424989, 0.01507831228824928, 0.01492204813964542, 0.01478643949493405, 0.01477233674211043, 0.0146562614272324, 0.01465227497604702)
$arr = New-Object 'object[,]' 1,$rowVals.Length
for ($i = 0; $i -lt $rowVals.Length; $i++) { $arr[0,$i] = $rowVals[$i] }
$ws.Range("E11:Y11").Value = $arr
